$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Simple value updates (rows whose coin identity is unchanged) ---
$ws.Range("D2").Value = "69.041.55"
$ws.Range("E2").Value = "  -0.12%  "
$ws.Range("D3").Value = "3.942.44"
$ws.Range("E3").Value = "  +3.53%  "
$ws.Range("E4").Value = "  -0.21%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "604.82"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "168.46"
$ws.Range("E6").Value = "  +2.53%  "
$ws.Range("D7").Value = "3.941.94"
$ws.Range("E7").Value = "  +3.55%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.535"
$ws.Range("E9").Value = "  +0.40%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.172"
$ws.Range("E10").Value = "  +2.37%  "
$ws.Range("E11").Value = "  +2.76%  "
$ws.Range("E12").Value = "  +1.54%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000258"
$ws.Range("E13").Value = "  +5.26%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "37.74"
$ws.Range("E14").Value = "  +1.93%  "
$ws.Range("D15").Value = "4.599.25"
$ws.Range("E15").Value = "  +3.52%  "
$ws.Range("D16").Value = "3.955.90"
$ws.Range("E16").Value = "  +3.83%  "
$ws.Range("D17").Value = "69.069.86"
$ws.Range("E17").Value = "  -0.29%  "
$ws.Range("E18").Value = "  +0.02%  "
$ws.Range("E19").Value = "  +1.42%  "
$ws.Range("E20").Value = "  -1.64%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.97"
$ws.Range("E21").Value = "  -3.76%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "495.94"
$ws.Range("E22").Value = "  +1.91%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.733"
$ws.Range("E23").Value = "  +1.84%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.0000170"
$ws.Range("E24").Value = "  +7.74%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "85.19"
$ws.Range("E25").Value = "  +0.88%  "
$ws.Range("E26").Value = "  +1.23%  "
$ws.Range("E27").Value = "  +0.55%  "
$ws.Range("E28").Value = "  +1.97%  "
$ws.Range("E29").Value = "  -0.01%  "
$ws.Range("E30").Value = "  +0.85%  "
$ws.Range("D31").Value = "4.093.32"
$ws.Range("E31").Value = "  +3.30%  "
$ws.Range("E32").Value = "  +0.36%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.83"
$ws.Range("E33").Value = "  -2.15%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "32.15"
$ws.Range("E34").Value = "  +1.15%  "
$ws.Range("D35").Value = "3.906.52"
$ws.Range("E35").Value = "  +4.16%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.108"
$ws.Range("E36").Value = "  +0.85%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.05"
$ws.Range("E37").Value = "  +2.26%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.02"
$ws.Range("E38").Value = "  +2.69%  "
$ws.Range("E41").Value = "  -0.08%  "
$ws.Range("E42").Value = "  +1.26%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "435.93"
$ws.Range("E43").Value = "  +0.31%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.01"
$ws.Range("E44").Value = "  +1.41%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "48.07"
$ws.Range("E45").Value = "  -1.02%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "8.63"
$ws.Range("E46").Value = "  +3.16%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "143.20"
$ws.Range("E49").Value = "  +0.29%  "
# --- Row swaps: ranking order changed between rows 39/40, 47/48, 50/51 ---
$ws.Range("B39").Value = "Kaspa"
$ws.Range("C39").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.139"
$ws.Range("E39").Value = "  -0.53%  "
$ws.Range("B40").Value = "dogwifhat"
$ws.Range("C40").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.30"
$ws.Range("E40").Value = "  +8.54%  "
$ws.Range("B47").Value = "FLOKI"
$ws.Range("C47").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.000288"
$ws.Range("E47").Value = "  +27.45%  "
$ws.Range("B48").Value = "USDe"
$ws.Range("C48").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.00"
$ws.Range("E48").Value = "  +0.02%  "
$ws.Range("B50").Value = "Maker"
$ws.Range("C50").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D50").Value = "2.820.51"
$ws.Range("E50").Value = "  -0.14%  "
$ws.Range("B51").Value = "VeChain"
$ws.Range("C51").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0360"
$ws.Range("E51").Value = "  +1.46%  "
